$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 45414.25
$ws.Range("J3").Value = 45414.25
$ws.Range("L3").Value = 45414.25
$ws.Range("N3").Value = -45642.25
$ws.Range("H82").Value = 1951.4286
$ws.Range("I82").Value = 1943.3334
$ws.Range("K82").Value = 5830.0002
$ws.Range("M82").Value = -5424.0002
$ws.Range("H85").Value = 1951.4286
$ws.Range("I85").Value = 1943.3334
$ws.Range("K85").Value = 5830.0002
$ws.Range("M85").Value = -4426.0002
$ws.Range("H86").Value = 200042130
$ws.Range("I86").Value = 250052300
$ws.Range("J86").Value = 1404
$ws.Range("K86").Value = 250052300
$ws.Range("L86").Value = 1404
$ws.Range("M86").Value = -250051177
$ws.Range("N86").Value = -3650
$ws.Range("H89").Value = 200042130
$ws.Range("I89").Value = 250052300
$ws.Range("J89").Value = 1404
$ws.Range("K89").Value = 1250261500
$ws.Range("L89").Value = 7020
$ws.Range("M89").Value = -1250255884
$ws.Range("N89").Value = -18252
$ws.Range("H102").Value = 45414.25
$ws.Range("J102").Value = 45414.25
$ws.Range("L102").Value = 45414.25
$ws.Range("N102").Value = -51904.25
$ws.Range("H121").Value = 1208.2812
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1208.2812
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3624.8436
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -7118.8436
$ws.Range("H133").Value = 66693.336
$ws.Range("J133").Value = 66693.336
$ws.Range("L133").Value = 66693.336
$ws.Range("N133").Value = -76813.336
$ws.Range("H137").Value = 4388483.5
$ws.Range("I137").Value = 6412339
$ws.Range("K137").Value = 19237017
$ws.Range("M137").Value = -19234467
$ws.Range("H138").Value = 4751.985
$ws.Range("I138").Value = 6211.75
$ws.Range("J138").Value = 4554.051
$ws.Range("K138").Value = 18635.25
$ws.Range("L138").Value = 13662.153
$ws.Range("M138").Value = -13495.25
$ws.Range("N138").Value = -23942.153

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 3395.8
$ws.Range("I35").Value = 3395.8
$ws.Range("K35").Value = 3395.8
$ws.Range("M35").Value = -2989.8
$ws.Range("H128").Value = 48357.25
$ws.Range("J128").Value = 48357.25
$ws.Range("L128").Value = 48357.25
$ws.Range("N128").Value = -58317.25
$ws.Range("H132").Value = 2202151.2
$ws.Range("I132").Value = 4670.4
$ws.Range("J132").Value = 5132125.5
$ws.Range("K132").Value = 14011.2
$ws.Range("L132").Value = 15396376.5
$ws.Range("M132").Value = -11481.2
$ws.Range("N132").Value = -15401436.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 978
$ws.Range("I12").Value = 978
$ws.Range("K12").Value = 978
$ws.Range("M12").Value = -810
$ws.Range("H134").Value = 2462.7666
$ws.Range("I134").Value = 2462.12
$ws.Range("J134").Value = 2466
$ws.Range("K134").Value = 7386.36
$ws.Range("L134").Value = 7398
$ws.Range("M134").Value = -4851.36
$ws.Range("N134").Value = -12468

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2754.2
$ws.Range("I16").Value = 2734.4167
$ws.Range("J16").Value = 2833.3333
$ws.Range("K16").Value = 2734.4167
$ws.Range("L16").Value = 2833.3333
$ws.Range("M16").Value = -2447.4167
$ws.Range("N16").Value = -3407.3333
$ws.Range("H31").Value = 11185.069
$ws.Range("I31").Value = 1500
$ws.Range("J31").Value = 11911.45
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 11911.45
$ws.Range("M31").Value = -1205
$ws.Range("N31").Value = -12501.45
$ws.Range("H34").Value = 11185.069
$ws.Range("I34").Value = 1500
$ws.Range("J34").Value = 11911.45
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 11911.45
$ws.Range("M34").Value = -1298
$ws.Range("N34").Value = -12315.45
$ws.Range("H113").Value = 2754.2
$ws.Range("I113").Value = 2734.4167
$ws.Range("J113").Value = 2833.3333
$ws.Range("K113").Value = 2734.4167
$ws.Range("L113").Value = 2833.3333
$ws.Range("M113").Value = -564.4167000000002
$ws.Range("N113").Value = -7173.3333
$ws.Range("H122").Value = 2163.6316
$ws.Range("I122").Value = 1515.8334
$ws.Range("K122").Value = 4547.5002
$ws.Range("M122").Value = -2097.5002
$ws.Range("H134").Value = 5438663.5
$ws.Range("I134").Value = 6101549
$ws.Range("J134").Value = 3002.8
$ws.Range("K134").Value = 18304647
$ws.Range("L134").Value = 9008.400000000001
$ws.Range("M134").Value = -18302112
$ws.Range("N134").Value = -14078.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 666.6667
$ws.Range("H36").Value = 2300
$ws.Range("I36").Value = 1500
$ws.Range("J36").Value = 2441.1765
$ws.Range("K36").Value = 4500
$ws.Range("L36").Value = 7323.529500000001
$ws.Range("M36").Value = -4331
$ws.Range("N36").Value = -7661.529500000001
$ws.Range("H55").Value = 1700.5333
$ws.Range("J55").Value = 1861.8462
$ws.Range("L55").Value = 5585.5386
$ws.Range("N55").Value = -5939.5386
$ws.Range("H113").Value = 919.4
$ws.Range("I113").Value = 744.2222
$ws.Range("J113").Value = 1182.1666
$ws.Range("K113").Value = 2232.6666
$ws.Range("L113").Value = 3546.4998
$ws.Range("M113").Value = -62.66660000000002
$ws.Range("N113").Value = -7886.4998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 170820
$ws.Range("I97").Value = 203984
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 203984
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -203488
$ws.Range("N97").Value = -5992
$ws.Range("H123").Value = 20518.572
$ws.Range("J123").Value = 20518.572
$ws.Range("L123").Value = 20518.572
$ws.Range("N123").Value = -25418.572

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H102").Value = 72112.2
$ws.Range("J102").Value = 72112.2
$ws.Range("L102").Value = 72112.2
$ws.Range("N102").Value = -78602.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 71221.75999999999
$ws.Range("J135").Value = 71221.75999999999
$ws.Range("L135").Value = 71221.75999999999
$ws.Range("N135").Value = -81361.75999999999
$ws.Range("H139").Value = 65074.285
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 65074.285
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 65074.285
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -75354.285
$ws.Range("H141").Value = 69723.81
$ws.Range("J141").Value = 69723.81
$ws.Range("L141").Value = 69723.81
$ws.Range("N141").Value = -80083.81
